$d = $word.ActiveDocument

# Pull the full package as flat OPC WordprocessingML (this round-trips the
# raw part contents, including word/numbering.xml, byte-for-byte when
# unmodified).
$xml = $d.WordOpenXML

# The numbering.xml part in this document contains two <w:abstractNum>
# elements that were both minted with w:abstractNumId="990" (one carrying
# nsid "170cd2de", the other "2c1ae401"), plus two identical <w:num
# w:numId="1000"> instances referencing abstractNumId 990. These are
# duplicate/clashing numbering definitions accidentally copied in from the
# reference template; remove the extra (first) copy of each so only one
# abstractNum and one num instance remain, matching what a non-duplicating
# writer would have produced.

$dupAbstractNum = '<w:abstractNum w:abstractNumId="990"><w:nsid w:val="170cd2de" /><w:multiLevelType w:val="multilevel" /><w:lvl w:ilvl="0"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="0" /></w:tabs><w:ind w:left="480" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="1"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="720" /></w:tabs><w:ind w:left="1200" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="2"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="1440" /></w:tabs><w:ind w:left="1920" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="3"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="2160" /></w:tabs><w:ind w:left="2640" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="4"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="2880" /></w:tabs><w:ind w:left="3360" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="5"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="3600" /></w:tabs><w:ind w:left="4080" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="6"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="4320" /></w:tabs><w:ind w:left="4800" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="7"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="5040" /></w:tabs><w:ind w:left="5520" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="8"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="5760" /></w:tabs><w:ind w:left="6240" w:hanging="480" /></w:pPr></w:lvl></w:abstractNum>'

if ($xml.IndexOf($dupAbstractNum) -ge 0) {
    $idx = $xml.IndexOf($dupAbstractNum)
    $xml = $xml.Substring(0, $idx) + $xml.Substring($idx + $dupAbstractNum.Length)
}

$dupNum = '<w:num w:numId="1000"><w:abstractNumId w:val="990" /></w:num>'
$bothNums = $dupNum + $dupNum
if ($xml.IndexOf($bothNums) -ge 0) {
    $idx = $xml.IndexOf($bothNums)
    $xml = $xml.Substring(0, $idx) + $dupNum + $xml.Substring($idx + $bothNums.Length)
}

$d.WordOpenXML = $xml
